$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns before column D (shifts existing D:K data to G:N)
$ws.Range("D5:F5").EntireColumn.Insert()

# Copy number formats from column G (the shifted original column D) into the new D:F
# columns, restricted to the row blocks that actually contain financial data.
$ws.Range("G7:G35").Copy()
$ws.Range("D7:F35").PasteSpecial(-4122)

$ws.Range("G38:G77").Copy()
$ws.Range("D38:F77").PasteSpecial(-4122)

$ws.Range("G80:G102").Copy()
$ws.Range("D80:F102").PasteSpecial(-4122)

# Row 7
$ws.Range("D7").Value2 = 43524
$ws.Range("E7").Value2 = 43434
$ws.Range("F7").Value2 = 43343
$ws.Range("G7").Value2 = 43251
$ws.Range("H7").Value2 = 43159
$ws.Range("I7").Value2 = 43069
$ws.Range("J7").Value2 = 42978
$ws.Range("K7").Value2 = 42886
$ws.Range("L7").Value2 = 42794
$ws.Range("M7").Value2 = 42704
$ws.Range("N7").Value2 = 42613
# Row 8
$ws.Range("D8").Value2 = 1402800
$ws.Range("E8").Value2 = 1277300
$ws.Range("F8").Value2 = 1308400
$ws.Range("G8").Value2 = 1204500
$ws.Range("H8").Value2 = 1054300
$ws.Range("I8").Value2 = 1076500
$ws.Range("J8").Value2 = 1084100
$ws.Range("K8").Value2 = 1044700
$ws.Range("L8").Value2 = 862200
$ws.Range("M8").Value2 = 994100
$ws.Range("N8").Value2 = 775600
# Row 9
$ws.Range("D9").Value2 = 1250100
$ws.Range("E9").Value2 = 1117600
$ws.Range("F9").Value2 = 1116000
$ws.Range("G9").Value2 = 1037400
$ws.Range("H9").Value2 = 927100
$ws.Range("I9").Value2 = 934100
$ws.Range("J9").Value2 = 977100
$ws.Range("K9").Value2 = 897900
$ws.Range("L9").Value2 = 725200
$ws.Range("M9").Value2 = 874900
$ws.Range("N9").Value2 = 647000
# Row 10
$ws.Range("D10").Value2 = 152700
$ws.Range("E10").Value2 = 159700
$ws.Range("F10").Value2 = 192400
$ws.Range("G10").Value2 = 167100
$ws.Range("H10").Value2 = 127200
$ws.Range("I10").Value2 = 142400
$ws.Range("J10").Value2 = 107000
$ws.Range("K10").Value2 = 146800
$ws.Range("L10").Value2 = 137000
$ws.Range("M10").Value2 = 119200
$ws.Range("N10").Value2 = 128600
# Row 12
$ws.Range("D12").Value2 = "NA"
$ws.Range("E12").Value2 = "NA"
$ws.Range("F12").Value2 = "NA"
$ws.Range("G12").Value2 = "NA"
$ws.Range("H12").Value2 = "NA"
$ws.Range("I12").Value2 = "NA"
$ws.Range("J12").Value2 = "NA"
$ws.Range("K12").Value2 = "NA"
$ws.Range("L12").Value2 = "NA"
$ws.Range("M12").Value2 = "NA"
$ws.Range("N12").Value2 = "NA"
# Row 13
$ws.Range("D13").Value2 = 0
$ws.Range("E13").Value2 = 0
$ws.Range("F13").Value2 = 0
$ws.Range("G13").Value2 = 0
$ws.Range("H13").Value2 = 0
$ws.Range("I13").Value2 = 0
$ws.Range("J13").Value2 = 0
$ws.Range("K13").Value2 = 0
$ws.Range("L13").Value2 = 0
$ws.Range("M13").Value2 = 0
$ws.Range("N13").Value2 = 0
# Row 14
$ws.Range("D14").Value2 = "NA"
$ws.Range("E14").Value2 = "NA"
$ws.Range("F14").Value2 = 1400
$ws.Range("G14").Value2 = 900
$ws.Range("H14").Value2 = 12000
$ws.Range("I14").Value2 = 100
$ws.Range("J14").Value2 = 24400
$ws.Range("K14").Value2 = 0
$ws.Range("L14").Value2 = "NA"
$ws.Range("M14").Value2 = "NA"
$ws.Range("N14").Value2 = 40000
# Row 15
$ws.Range("D15").Value2 = 0
$ws.Range("E15").Value2 = 0
$ws.Range("F15").Value2 = 0
$ws.Range("G15").Value2 = 0
$ws.Range("H15").Value2 = 0
$ws.Range("I15").Value2 = 0
$ws.Range("J15").Value2 = 0
$ws.Range("K15").Value2 = 0
$ws.Range("L15").Value2 = 0
$ws.Range("M15").Value2 = 0
$ws.Range("N15").Value2 = 0
# Row 17
$ws.Range("D17").Value2 = 1348300
$ws.Range("E17").Value2 = 1234900
$ws.Range("F17").Value2 = 1225800
$ws.Range("G17").Value2 = 1139400
$ws.Range("H17").Value2 = 1033800
$ws.Range("I17").Value2 = 1032600
$ws.Range("J17").Value2 = 1106500
$ws.Range("K17").Value2 = 990300
$ws.Range("L17").Value2 = 819900
$ws.Range("M17").Value2 = 974400
$ws.Range("N17").Value2 = 788600
# Row 18
$ws.Range("D18").Value2 = 54500
$ws.Range("E18").Value2 = 42400
$ws.Range("F18").Value2 = 82600
$ws.Range("G18").Value2 = 65100
$ws.Range("H18").Value2 = 20500
$ws.Range("I18").Value2 = 43900
$ws.Range("J18").Value2 = -22400
$ws.Range("K18").Value2 = 54400
$ws.Range("L18").Value2 = 42300
$ws.Range("M18").Value2 = 19700
$ws.Range("N18").Value2 = -13000
# Row 20
$ws.Range("D20").Value2 = -3000
$ws.Range("E20").Value2 = -700
$ws.Range("F20").Value2 = -9000
$ws.Range("G20").Value2 = 2000
$ws.Range("H20").Value2 = -1800
$ws.Range("I20").Value2 = 3000
$ws.Range("J20").Value2 = 12300
$ws.Range("K20").Value2 = 600
$ws.Range("L20").Value2 = 900
$ws.Range("M20").Value2 = 600
$ws.Range("N20").Value2 = -1900
# Row 21
$ws.Range("D21").Value2 = 92800
$ws.Range("E21").Value2 = 76900
$ws.Range("F21").Value2 = 105800
$ws.Range("G21").Value2 = 100300
$ws.Range("H21").Value2 = 52800
$ws.Range("I21").Value2 = 79100
$ws.Range("J21").Value2 = 21900
$ws.Range("K21").Value2 = 87300
$ws.Range("L21").Value2 = 73700
$ws.Range("M21").Value2 = 50600
$ws.Range("N21").Value2 = 16600
# Row 22
$ws.Range("D22").Value2 = 18500
$ws.Range("E22").Value2 = 16700
$ws.Range("F22").Value2 = 15700
$ws.Range("G22").Value2 = 11500
$ws.Range("H22").Value2 = 7200
$ws.Range("I22").Value2 = 6600
$ws.Range("J22").Value2 = 5900
$ws.Range("K22").Value2 = 12400
$ws.Range("L22").Value2 = 12400
$ws.Range("M22").Value2 = 13300
$ws.Range("N22").Value2 = 12500
# Row 23
$ws.Range("D23").Value2 = 33100
$ws.Range("E23").Value2 = 25000
$ws.Range("F23").Value2 = 57900
$ws.Range("G23").Value2 = 55600
$ws.Range("H23").Value2 = 11500
$ws.Range("I23").Value2 = 40300
$ws.Range("J23").Value2 = -16000
$ws.Range("K23").Value2 = 42600
$ws.Range("L23").Value2 = 30800
$ws.Range("M23").Value2 = 7000
$ws.Range("N23").Value2 = -27400
# Row 24
$ws.Range("D24").Value2 = 10500
$ws.Range("E24").Value2 = 5600
$ws.Range("F24").Value2 = 5600
$ws.Range("G24").Value2 = 13300
$ws.Range("H24").Value2 = -8900
$ws.Range("I24").Value2 = 8400
$ws.Range("J24").Value2 = -6000
$ws.Range("K24").Value2 = 11000
$ws.Range("L24").Value2 = 7800
$ws.Range("M24").Value2 = 2100
$ws.Range("N24").Value2 = -13700
# Row 25
$ws.Range("D25").Value2 = 0
$ws.Range("E25").Value2 = 0
$ws.Range("F25").Value2 = 0
$ws.Range("G25").Value2 = 0
$ws.Range("H25").Value2 = 0
$ws.Range("I25").Value2 = 0
$ws.Range("J25").Value2 = 0
$ws.Range("K25").Value2 = 0
$ws.Range("L25").Value2 = 0
$ws.Range("M25").Value2 = 0
$ws.Range("N25").Value2 = 0
# Row 26
$ws.Range("D26").Value2 = 22500
$ws.Range("E26").Value2 = 19400
$ws.Range("F26").Value2 = 52400
$ws.Range("G26").Value2 = 42300
$ws.Range("H26").Value2 = 20400
$ws.Range("I26").Value2 = 31900
$ws.Range("J26").Value2 = -10100
$ws.Range("K26").Value2 = 31600
$ws.Range("L26").Value2 = 23000
$ws.Range("M26").Value2 = 4900
$ws.Range("N26").Value2 = -13700
# Row 27
$ws.Range("D27").Value2 = 22500
$ws.Range("E27").Value2 = 19400
$ws.Range("F27").Value2 = 52400
$ws.Range("G27").Value2 = 42300
$ws.Range("H27").Value2 = 20400
$ws.Range("I27").Value2 = 31900
$ws.Range("J27").Value2 = -10100
$ws.Range("K27").Value2 = 31600
$ws.Range("L27").Value2 = 23000
$ws.Range("M27").Value2 = 4900
$ws.Range("N27").Value2 = -13700
# Row 28
$ws.Range("D28").Value2 = 0
$ws.Range("E28").Value2 = 0
$ws.Range("F28").Value2 = 0
$ws.Range("G28").Value2 = 0
$ws.Range("H28").Value2 = 0
$ws.Range("I28").Value2 = 0
$ws.Range("J28").Value2 = 0
$ws.Range("K28").Value2 = 0
$ws.Range("L28").Value2 = 0
$ws.Range("M28").Value2 = 0
$ws.Range("N28").Value2 = 0
# Row 29
$ws.Range("D29").Value2 = -8700
$ws.Range("E29").Value2 = 300
$ws.Range("F29").Value2 = -800
$ws.Range("G29").Value2 = -2400
$ws.Range("H29").Value2 = -10200
$ws.Range("I29").Value2 = 4900
$ws.Range("J29").Value2 = -19500
$ws.Range("K29").Value2 = 7700
$ws.Range("L29").Value2 = 7300
$ws.Range("M29").Value2 = 1300
$ws.Range("N29").Value2 = 13600
# Row 30
$ws.Range("D30").Value2 = 0
$ws.Range("E30").Value2 = 0
$ws.Range("F30").Value2 = 0
$ws.Range("G30").Value2 = 0
$ws.Range("H30").Value2 = 0
$ws.Range("I30").Value2 = 0
$ws.Range("J30").Value2 = 0
$ws.Range("K30").Value2 = 0
$ws.Range("L30").Value2 = 0
$ws.Range("M30").Value2 = 0
$ws.Range("N30").Value2 = 0
# Row 31
$ws.Range("D31").Value2 = 0
$ws.Range("E31").Value2 = 0
$ws.Range("F31").Value2 = 0
$ws.Range("G31").Value2 = 0
$ws.Range("H31").Value2 = 0
$ws.Range("I31").Value2 = 0
$ws.Range("J31").Value2 = 0
$ws.Range("K31").Value2 = 0
$ws.Range("L31").Value2 = 0
$ws.Range("M31").Value2 = 0
$ws.Range("N31").Value2 = 0
# Row 32
$ws.Range("D32").Value2 = 3000
$ws.Range("E32").Value2 = 700
$ws.Range("F32").Value2 = 9000
$ws.Range("G32").Value2 = -2000
$ws.Range("H32").Value2 = 1800
$ws.Range("I32").Value2 = -3000
$ws.Range("J32").Value2 = -12300
$ws.Range("K32").Value2 = -600
$ws.Range("L32").Value2 = -900
$ws.Range("M32").Value2 = -600
$ws.Range("N32").Value2 = 1900
# Row 33
$ws.Range("D33").Value2 = 13900
$ws.Range("E33").Value2 = 19700
$ws.Range("F33").Value2 = 51600
$ws.Range("G33").Value2 = 40000
$ws.Range("H33").Value2 = 10200
$ws.Range("I33").Value2 = 36800
$ws.Range("J33").Value2 = -29500
$ws.Range("K33").Value2 = 39300
$ws.Range("L33").Value2 = 30300
$ws.Range("M33").Value2 = 6300
$ws.Range("N33").Value2 = -100
# Row 34
$ws.Range("D34").Value2 = 0
$ws.Range("E34").Value2 = 0
$ws.Range("F34").Value2 = 0
$ws.Range("G34").Value2 = 0
$ws.Range("H34").Value2 = 0
$ws.Range("I34").Value2 = 0
$ws.Range("J34").Value2 = 0
$ws.Range("K34").Value2 = 0
$ws.Range("L34").Value2 = 0
$ws.Range("M34").Value2 = 0
$ws.Range("N34").Value2 = 0
# Row 35
$ws.Range("D35").Value2 = 13900
$ws.Range("E35").Value2 = 19700
$ws.Range("F35").Value2 = 51600
$ws.Range("G35").Value2 = 40000
$ws.Range("H35").Value2 = 10200
$ws.Range("I35").Value2 = 36800
$ws.Range("J35").Value2 = -29500
$ws.Range("K35").Value2 = 39300
$ws.Range("L35").Value2 = 30300
$ws.Range("M35").Value2 = 6300
$ws.Range("N35").Value2 = -100
# Row 38
$ws.Range("D38").Value2 = 43524
$ws.Range("E38").Value2 = 43434
$ws.Range("F38").Value2 = 43343
$ws.Range("G38").Value2 = 43251
$ws.Range("H38").Value2 = 43159
$ws.Range("I38").Value2 = 43069
$ws.Range("J38").Value2 = 42978
$ws.Range("K38").Value2 = 42886
$ws.Range("L38").Value2 = 42794
$ws.Range("M38").Value2 = 42704
$ws.Range("N38").Value2 = 42613
# Row 41
$ws.Range("D41").Value2 = 66700
$ws.Range("E41").Value2 = 52400
$ws.Range("F41").Value2 = 622500
$ws.Range("G41").Value2 = 600400
$ws.Range("H41").Value2 = 195200
$ws.Range("I41").Value2 = 130200
$ws.Range("J41").Value2 = 252600
$ws.Range("K41").Value2 = 275800
$ws.Range("L41").Value2 = 395500
$ws.Range("M41").Value2 = 465200
$ws.Range("N41").Value2 = 517500
# Row 42
$ws.Range("D42").Value2 = 0
$ws.Range("E42").Value2 = 0
$ws.Range("F42").Value2 = 0
$ws.Range("G42").Value2 = 0
$ws.Range("H42").Value2 = 0
$ws.Range("I42").Value2 = 0
$ws.Range("J42").Value2 = 0
$ws.Range("K42").Value2 = "NA"
$ws.Range("L42").Value2 = 200
$ws.Range("M42").Value2 = 100
$ws.Range("N42").Value2 = 0
# Row 43
$ws.Range("D43").Value2 = 1053000
$ws.Range("E43").Value2 = 1088500
$ws.Range("F43").Value2 = 749500
$ws.Range("G43").Value2 = 678300
$ws.Range("H43").Value2 = 634700
$ws.Range("I43").Value2 = 772600
$ws.Range("J43").Value2 = 561400
$ws.Range("K43").Value2 = 870000
$ws.Range("L43").Value2 = 774300
$ws.Range("M43").Value2 = 716600
$ws.Range("N43").Value2 = 765800
# Row 44
$ws.Range("D44").Value2 = 866400
$ws.Range("E44").Value2 = 828600
$ws.Range("F44").Value2 = 589000
$ws.Range("G44").Value2 = 595200
$ws.Range("H44").Value2 = 523400
$ws.Range("I44").Value2 = 564800
$ws.Range("J44").Value2 = 462600
$ws.Range("K44").Value2 = 798000
$ws.Range("L44").Value2 = 720800
$ws.Range("M44").Value2 = 633800
$ws.Range("N44").Value2 = 652800
# Row 45
$ws.Range("D45").Value2 = 84100
$ws.Range("E45").Value2 = 59600
$ws.Range("F45").Value2 = 116200
$ws.Range("G45").Value2 = 120900
$ws.Range("H45").Value2 = 294700
$ws.Range("I45").Value2 = 132900
$ws.Range("J45").Value2 = 437200
$ws.Range("K45").Value2 = 108200
$ws.Range("L45").Value2 = 96300
$ws.Range("M45").Value2 = 97000
$ws.Range("N45").Value2 = 112000
# Row 46
$ws.Range("D46").Value2 = 2070300
$ws.Range("E46").Value2 = 2028900
$ws.Range("F46").Value2 = 2077200
$ws.Range("G46").Value2 = 1995000
$ws.Range("H46").Value2 = 1648000
$ws.Range("I46").Value2 = 1600500
$ws.Range("J46").Value2 = 1713900
$ws.Range("K46").Value2 = 2052000
$ws.Range("L46").Value2 = 1987000
$ws.Range("M46").Value2 = 1912700
$ws.Range("N46").Value2 = 2048100
# Row 47
$ws.Range("D47").Value2 = 0
$ws.Range("E47").Value2 = 0
$ws.Range("F47").Value2 = 0
$ws.Range("G47").Value2 = 0
$ws.Range("H47").Value2 = 0
$ws.Range("I47").Value2 = 0
$ws.Range("J47").Value2 = 0
$ws.Range("K47").Value2 = 0
$ws.Range("L47").Value2 = 0
$ws.Range("M47").Value2 = 0
$ws.Range("N47").Value2 = 0
# Row 48
$ws.Range("D48").Value2 = 1478300
$ws.Range("E48").Value2 = 1492200
$ws.Range("F48").Value2 = 1075000
$ws.Range("G48").Value2 = 1074400
$ws.Range("H48").Value2 = 1083200
$ws.Range("I48").Value2 = 1092800
$ws.Range("J48").Value2 = 1051700
$ws.Range("K48").Value2 = 1016900
$ws.Range("L48").Value2 = 940300
$ws.Range("M48").Value2 = 893200
$ws.Range("N48").Value2 = 895000
# Row 49
$ws.Range("D49").Value2 = 64300
$ws.Range("E49").Value2 = 64300
$ws.Range("F49").Value2 = 64300
$ws.Range("G49").Value2 = 64300
$ws.Range("H49").Value2 = 64500
$ws.Range("I49").Value2 = 64900
$ws.Range("J49").Value2 = 64900
$ws.Range("K49").Value2 = 66800
$ws.Range("L49").Value2 = 66500
$ws.Range("M49").Value2 = 66100
$ws.Range("N49").Value2 = 66400
# Row 50
$ws.Range("D50").Value2 = 0
$ws.Range("E50").Value2 = 0
$ws.Range("F50").Value2 = 0
$ws.Range("G50").Value2 = 0
$ws.Range("H50").Value2 = 0
$ws.Range("I50").Value2 = 0
$ws.Range("J50").Value2 = 0
$ws.Range("K50").Value2 = 0
$ws.Range("L50").Value2 = 0
$ws.Range("M50").Value2 = 0
$ws.Range("N50").Value2 = 0
# Row 51
$ws.Range("D51").Value2 = 0
$ws.Range("E51").Value2 = 0
$ws.Range("F51").Value2 = 0
$ws.Range("G51").Value2 = 0
$ws.Range("H51").Value2 = 0
$ws.Range("I51").Value2 = 0
$ws.Range("J51").Value2 = 0
$ws.Range("K51").Value2 = 0
$ws.Range("L51").Value2 = 0
$ws.Range("M51").Value2 = 0
$ws.Range("N51").Value2 = 0
# Row 52
$ws.Range("D52").Value2 = 115900
$ws.Range("E52").Value2 = 123200
$ws.Range("F52").Value2 = 111800
$ws.Range("G52").Value2 = 111900
$ws.Range("H52").Value2 = 114700
$ws.Range("I52").Value2 = 140300
$ws.Range("J52").Value2 = 144600
$ws.Range("K52").Value2 = 139000
$ws.Range("L52").Value2 = 137900
$ws.Range("M52").Value2 = 130600
$ws.Range("N52").Value2 = 121300
# Row 53
$ws.Range("D53").Value2 = 0
$ws.Range("E53").Value2 = 0
$ws.Range("F53").Value2 = 0
$ws.Range("G53").Value2 = 0
$ws.Range("H53").Value2 = 0
$ws.Range("I53").Value2 = 0
$ws.Range("J53").Value2 = 0
$ws.Range("K53").Value2 = 0
$ws.Range("L53").Value2 = 0
$ws.Range("M53").Value2 = 0
$ws.Range("N53").Value2 = 0
# Row 54
$ws.Range("D54").Value2 = 3728700
$ws.Range("E54").Value2 = 3708700
$ws.Range("F54").Value2 = 3328300
$ws.Range("G54").Value2 = 3245500
$ws.Range("H54").Value2 = 2910500
$ws.Range("I54").Value2 = 2898600
$ws.Range("J54").Value2 = 2975100
$ws.Range("K54").Value2 = 3274600
$ws.Range("L54").Value2 = 3131800
$ws.Range("M54").Value2 = 3002600
$ws.Range("N54").Value2 = 3130900
# Row 57
$ws.Range("D57").Value2 = 322100
$ws.Range("E57").Value2 = 319600
$ws.Range("F57").Value2 = 261300
$ws.Range("G57").Value2 = 241600
$ws.Range("H57").Value2 = 247600
$ws.Range("I57").Value2 = 244700
$ws.Range("J57").Value2 = 226500
$ws.Range("K57").Value2 = 346500
$ws.Range("L57").Value2 = 307500
$ws.Range("M57").Value2 = 224700
$ws.Range("N57").Value2 = 243500
# Row 58
$ws.Range("D58").Value2 = 88900
$ws.Range("E58").Value2 = 29100
$ws.Range("F58").Value2 = 19700
$ws.Range("G58").Value2 = 19900
$ws.Range("H58").Value2 = 19000
$ws.Range("I58").Value2 = 21200
$ws.Range("J58").Value2 = 19200
$ws.Range("K58").Value2 = 311700
$ws.Range("L58").Value2 = 312200
$ws.Range("M58").Value2 = 312900
$ws.Range("N58").Value2 = 313500
# Row 59
$ws.Range("D59").Value2 = 341300
$ws.Range("E59").Value2 = 409700
$ws.Range("F59").Value2 = 260900
$ws.Range("G59").Value2 = 250500
$ws.Range("H59").Value2 = 263800
$ws.Range("I59").Value2 = 234200
$ws.Range("J59").Value2 = 362800
$ws.Range("K59").Value2 = 258300
$ws.Range("L59").Value2 = 220400
$ws.Range("M59").Value2 = 202800
$ws.Range("N59").Value2 = 264100
# Row 60
$ws.Range("D60").Value2 = 752300
$ws.Range("E60").Value2 = 758500
$ws.Range("F60").Value2 = 541900
$ws.Range("G60").Value2 = 511900
$ws.Range("H60").Value2 = 530300
$ws.Range("I60").Value2 = 500100
$ws.Range("J60").Value2 = 608400
$ws.Range("K60").Value2 = 916500
$ws.Range("L60").Value2 = 840100
$ws.Range("M60").Value2 = 740500
$ws.Range("N60").Value2 = 821100
# Row 61
$ws.Range("D61").Value2 = 1310200
$ws.Range("E61").Value2 = 1307800
$ws.Range("F61").Value2 = 1138600
$ws.Range("G61").Value2 = 1139100
$ws.Range("H61").Value2 = 799800
$ws.Range("I61").Value2 = 803800
$ws.Range("J61").Value2 = 805600
$ws.Range("K61").Value2 = 751700
$ws.Range("L61").Value2 = 752100
$ws.Range("M61").Value2 = 755200
$ws.Range("N61").Value2 = 757900
# Row 62
$ws.Range("D62").Value2 = 167700
$ws.Range("E62").Value2 = 153300
$ws.Range("F62").Value2 = 154200
$ws.Range("G62").Value2 = 141600
$ws.Range("H62").Value2 = 128800
$ws.Range("I62").Value2 = 159700
$ws.Range("J62").Value2 = 160200
$ws.Range("K62").Value2 = 188400
$ws.Range("L62").Value2 = 177600
$ws.Range("M62").Value2 = 176900
$ws.Range("N62").Value2 = 184400
# Row 63
$ws.Range("D63").Value2 = 0
$ws.Range("E63").Value2 = 0
$ws.Range("F63").Value2 = 0
$ws.Range("G63").Value2 = 0
$ws.Range("H63").Value2 = 0
$ws.Range("I63").Value2 = 0
$ws.Range("J63").Value2 = 0
$ws.Range("K63").Value2 = 0
$ws.Range("L63").Value2 = 0
$ws.Range("M63").Value2 = 0
$ws.Range("N63").Value2 = 0
# Row 64
$ws.Range("D64").Value2 = 0
$ws.Range("E64").Value2 = 0
$ws.Range("F64").Value2 = 0
$ws.Range("G64").Value2 = 0
$ws.Range("H64").Value2 = 0
$ws.Range("I64").Value2 = 0
$ws.Range("J64").Value2 = 0
$ws.Range("K64").Value2 = 0
$ws.Range("L64").Value2 = 0
$ws.Range("M64").Value2 = 0
$ws.Range("N64").Value2 = 0
# Row 65
$ws.Range("D65").Value2 = 0
$ws.Range("E65").Value2 = 0
$ws.Range("F65").Value2 = 0
$ws.Range("G65").Value2 = 0
$ws.Range("H65").Value2 = 0
$ws.Range("I65").Value2 = 0
$ws.Range("J65").Value2 = 0
$ws.Range("K65").Value2 = 0
$ws.Range("L65").Value2 = 0
$ws.Range("M65").Value2 = 0
$ws.Range("N65").Value2 = 0
# Row 66
$ws.Range("D66").Value2 = 2230400
$ws.Range("E66").Value2 = 2219800
$ws.Range("F66").Value2 = 1834900
$ws.Range("G66").Value2 = 1792800
$ws.Range("H66").Value2 = 1459200
$ws.Range("I66").Value2 = 1463800
$ws.Range("J66").Value2 = 1574400
$ws.Range("K66").Value2 = 1856700
$ws.Range("L66").Value2 = 1770000
$ws.Range("M66").Value2 = 1672600
$ws.Range("N66").Value2 = 1763600
# Row 68
$ws.Range("D68").Value2 = 0
$ws.Range("E68").Value2 = 0
$ws.Range("F68").Value2 = 0
$ws.Range("G68").Value2 = 0
$ws.Range("H68").Value2 = 0
$ws.Range("I68").Value2 = 0
$ws.Range("J68").Value2 = 0
$ws.Range("K68").Value2 = 0
$ws.Range("L68").Value2 = 0
$ws.Range("M68").Value2 = 0
$ws.Range("N68").Value2 = 0
# Row 69
$ws.Range("D69").Value2 = 0
$ws.Range("E69").Value2 = 0
$ws.Range("F69").Value2 = 0
$ws.Range("G69").Value2 = 0
$ws.Range("H69").Value2 = 0
$ws.Range("I69").Value2 = 0
$ws.Range("J69").Value2 = 0
$ws.Range("K69").Value2 = 0
$ws.Range("L69").Value2 = 0
$ws.Range("M69").Value2 = 0
$ws.Range("N69").Value2 = 0
# Row 70
$ws.Range("D70").Value2 = 0
$ws.Range("E70").Value2 = 0
$ws.Range("F70").Value2 = 0
$ws.Range("G70").Value2 = 0
$ws.Range("H70").Value2 = 0
$ws.Range("I70").Value2 = 0
$ws.Range("J70").Value2 = 0
$ws.Range("K70").Value2 = 0
$ws.Range("L70").Value2 = 0
$ws.Range("M70").Value2 = 0
$ws.Range("N70").Value2 = 0
# Row 71
$ws.Range("D71").Value2 = 0
$ws.Range("E71").Value2 = 0
$ws.Range("F71").Value2 = 0
$ws.Range("G71").Value2 = 0
$ws.Range("H71").Value2 = 0
$ws.Range("I71").Value2 = 0
$ws.Range("J71").Value2 = 0
$ws.Range("K71").Value2 = 0
$ws.Range("L71").Value2 = 0
$ws.Range("M71").Value2 = 0
$ws.Range("N71").Value2 = 0
# Row 72
$ws.Range("D72").Value2 = 1449200
$ws.Range("E72").Value2 = 1449400
$ws.Range("F72").Value2 = 1446500
$ws.Range("G72").Value2 = 1408700
$ws.Range("H72").Value2 = 1382800
$ws.Range("I72").Value2 = 1386600
$ws.Range("J72").Value2 = 1363800
$ws.Range("K72").Value2 = 1407200
$ws.Range("L72").Value2 = 1381900
$ws.Range("M72").Value2 = 1365400
$ws.Range("N72").Value2 = 1373000
# Row 73
$ws.Range("D73").Value2 = 0
$ws.Range("E73").Value2 = 0
$ws.Range("F73").Value2 = 0
$ws.Range("G73").Value2 = 0
$ws.Range("H73").Value2 = 0
$ws.Range("I73").Value2 = 0
$ws.Range("J73").Value2 = 0
$ws.Range("K73").Value2 = 0
$ws.Range("L73").Value2 = 0
$ws.Range("M73").Value2 = 0
$ws.Range("N73").Value2 = 0
# Row 74
$ws.Range("D74").Value2 = 0
$ws.Range("E74").Value2 = 0
$ws.Range("F74").Value2 = 0
$ws.Range("G74").Value2 = 0
$ws.Range("H74").Value2 = 0
$ws.Range("I74").Value2 = 0
$ws.Range("J74").Value2 = 0
$ws.Range("K74").Value2 = 0
$ws.Range("L74").Value2 = 0
$ws.Range("M74").Value2 = 0
$ws.Range("N74").Value2 = 0
# Row 75
$ws.Range("D75").Value2 = 0
$ws.Range("E75").Value2 = 0
$ws.Range("F75").Value2 = 0
$ws.Range("G75").Value2 = 0
$ws.Range("H75").Value2 = 0
$ws.Range("I75").Value2 = 0
$ws.Range("J75").Value2 = 0
$ws.Range("K75").Value2 = 0
$ws.Range("L75").Value2 = 0
$ws.Range("M75").Value2 = 0
$ws.Range("N75").Value2 = 0
# Row 76
$ws.Range("D76").Value2 = 1498300
$ws.Range("E76").Value2 = 1488800
$ws.Range("F76").Value2 = 1493400
$ws.Range("G76").Value2 = 1452700
$ws.Range("H76").Value2 = 1451300
$ws.Range("I76").Value2 = 1434800
$ws.Range("J76").Value2 = 1400800
$ws.Range("K76").Value2 = 1417900
$ws.Range("L76").Value2 = 1361800
$ws.Range("M76").Value2 = 1329900
$ws.Range("N76").Value2 = 1367300
# Row 77
$ws.Range("D77").Value2 = 0
$ws.Range("E77").Value2 = 0
$ws.Range("F77").Value2 = 0
$ws.Range("G77").Value2 = 0
$ws.Range("H77").Value2 = 0
$ws.Range("I77").Value2 = 0
$ws.Range("J77").Value2 = 0
$ws.Range("K77").Value2 = 0
$ws.Range("L77").Value2 = 0
$ws.Range("M77").Value2 = 0
$ws.Range("N77").Value2 = 0
# Row 80
$ws.Range("D80").Value2 = 43524
$ws.Range("E80").Value2 = 43434
$ws.Range("F80").Value2 = 43343
$ws.Range("G80").Value2 = 43251
$ws.Range("H80").Value2 = 43159
$ws.Range("I80").Value2 = 43069
$ws.Range("J80").Value2 = 42978
$ws.Range("K80").Value2 = 42886
$ws.Range("L80").Value2 = 42794
$ws.Range("M80").Value2 = 42704
$ws.Range("N80").Value2 = 42613
# Row 81
$ws.Range("D81").Value2 = 13900
$ws.Range("E81").Value2 = 19700
$ws.Range("F81").Value2 = 51600
$ws.Range("G81").Value2 = 40000
$ws.Range("H81").Value2 = 10200
$ws.Range("I81").Value2 = 36800
$ws.Range("J81").Value2 = -29500
$ws.Range("K81").Value2 = 39300
$ws.Range("L81").Value2 = 30300
$ws.Range("M81").Value2 = 6300
$ws.Range("N81").Value2 = -100
# Row 83
$ws.Range("D83").Value2 = 41200
$ws.Range("E83").Value2 = 35200
$ws.Range("F83").Value2 = 32200
$ws.Range("G83").Value2 = 33100
$ws.Range("H83").Value2 = 34100
$ws.Range("I83").Value2 = 32200
$ws.Range("J83").Value2 = 32000
$ws.Range("K83").Value2 = 32300
$ws.Range("L83").Value2 = 30500
$ws.Range("M83").Value2 = 30300
$ws.Range("N83").Value2 = 31500
# Row 84
$ws.Range("D84").Value2 = 0
$ws.Range("E84").Value2 = 0
$ws.Range("F84").Value2 = 0
$ws.Range("G84").Value2 = 0
$ws.Range("H84").Value2 = 0
$ws.Range("I84").Value2 = 0
$ws.Range("J84").Value2 = 0
$ws.Range("K84").Value2 = 0
$ws.Range("L84").Value2 = 0
$ws.Range("M84").Value2 = 0
$ws.Range("N84").Value2 = 0
# Row 85
$ws.Range("D85").Value2 = 0
$ws.Range("E85").Value2 = 0
$ws.Range("F85").Value2 = 0
$ws.Range("G85").Value2 = 0
$ws.Range("H85").Value2 = 0
$ws.Range("I85").Value2 = 0
$ws.Range("J85").Value2 = 0
$ws.Range("K85").Value2 = 0
$ws.Range("L85").Value2 = 0
$ws.Range("M85").Value2 = 0
$ws.Range("N85").Value2 = 0
# Row 86
$ws.Range("D86").Value2 = 0
$ws.Range("E86").Value2 = 0
$ws.Range("F86").Value2 = 0
$ws.Range("G86").Value2 = 0
$ws.Range("H86").Value2 = 0
$ws.Range("I86").Value2 = 0
$ws.Range("J86").Value2 = 0
$ws.Range("K86").Value2 = 0
$ws.Range("L86").Value2 = 0
$ws.Range("M86").Value2 = 0
$ws.Range("N86").Value2 = 0
# Row 87
$ws.Range("D87").Value2 = 0
$ws.Range("E87").Value2 = 0
$ws.Range("F87").Value2 = 0
$ws.Range("G87").Value2 = 0
$ws.Range("H87").Value2 = 0
$ws.Range("I87").Value2 = 0
$ws.Range("J87").Value2 = 0
$ws.Range("K87").Value2 = 0
$ws.Range("L87").Value2 = 0
$ws.Range("M87").Value2 = 0
$ws.Range("N87").Value2 = 0
# Row 88
$ws.Range("D88").Value2 = 0
$ws.Range("E88").Value2 = 0
$ws.Range("F88").Value2 = 0
$ws.Range("G88").Value2 = 0
$ws.Range("H88").Value2 = 0
$ws.Range("I88").Value2 = 0
$ws.Range("J88").Value2 = 0
$ws.Range("K88").Value2 = 0
$ws.Range("L88").Value2 = 0
$ws.Range("M88").Value2 = 0
$ws.Range("N88").Value2 = 0
# Row 89
$ws.Range("D89").Value2 = 4700
$ws.Range("E89").Value2 = -357600
$ws.Range("F89").Value2 = 67000
$ws.Range("G89").Value2 = 42900
$ws.Range("H89").Value2 = 83000
$ws.Range("I89").Value2 = -34200
$ws.Range("J89").Value2 = 155600
$ws.Range("K89").Value2 = 6900
$ws.Range("L89").Value2 = 12900
$ws.Range("M89").Value2 = -1000
$ws.Range("N89").Value2 = 80000
# Row 91
$ws.Range("D91").Value2 = -29600
$ws.Range("E91").Value2 = -37900
$ws.Range("F91").Value2 = -30400
$ws.Range("G91").Value2 = -43200
$ws.Range("H91").Value2 = -41300
$ws.Range("I91").Value2 = -59700
$ws.Range("J91").Value2 = -51000
$ws.Range("K91").Value2 = -71300
$ws.Range("L91").Value2 = -47800
$ws.Range("M91").Value2 = -43000
$ws.Range("N91").Value2 = -58900
# Row 92
$ws.Range("D92").Value2 = 0
$ws.Range("E92").Value2 = 0
$ws.Range("F92").Value2 = 0
$ws.Range("G92").Value2 = 0
$ws.Range("H92").Value2 = 0
$ws.Range("I92").Value2 = 0
$ws.Range("J92").Value2 = 0
$ws.Range("K92").Value2 = 0
$ws.Range("L92").Value2 = 0
$ws.Range("M92").Value2 = 0
$ws.Range("N92").Value2 = 0
# Row 93
$ws.Range("D93").Value2 = 0
$ws.Range("E93").Value2 = 0
$ws.Range("F93").Value2 = 0
$ws.Range("G93").Value2 = 0
$ws.Range("H93").Value2 = 0
$ws.Range("I93").Value2 = 0
$ws.Range("J93").Value2 = 0
$ws.Range("K93").Value2 = 0
$ws.Range("L93").Value2 = 0
$ws.Range("M93").Value2 = 0
$ws.Range("N93").Value2 = 0
# Row 94
$ws.Range("D94").Value2 = -35700
$ws.Range("E94").Value2 = -357400
$ws.Range("F94").Value2 = -49800
$ws.Range("G94").Value2 = -178300
$ws.Range("H94").Value2 = 135300
$ws.Range("I94").Value2 = 22100
$ws.Range("J94").Value2 = 104500
$ws.Range("K94").Value2 = -113200
$ws.Range("L94").Value2 = -68300
$ws.Range("M94").Value2 = -25700
$ws.Range("N94").Value2 = -25500
# Row 96
$ws.Range("D96").Value2 = -14100
$ws.Range("E96").Value2 = -14100
$ws.Range("F96").Value2 = -14000
$ws.Range("G96").Value2 = -14000
$ws.Range("H96").Value2 = -14000
$ws.Range("I96").Value2 = -14000
$ws.Range("J96").Value2 = -13900
$ws.Range("K96").Value2 = -13900
$ws.Range("L96").Value2 = -13900
$ws.Range("M96").Value2 = -13900
$ws.Range("N96").Value2 = -13800
# Row 97
$ws.Range("D97").Value2 = 0
$ws.Range("E97").Value2 = 0
$ws.Range("F97").Value2 = 0
$ws.Range("G97").Value2 = 0
$ws.Range("H97").Value2 = 0
$ws.Range("I97").Value2 = 0
$ws.Range("J97").Value2 = 0
$ws.Range("K97").Value2 = 0
$ws.Range("L97").Value2 = 0
$ws.Range("M97").Value2 = 0
$ws.Range("N97").Value2 = 0
# Row 98
$ws.Range("D98").Value2 = 0
$ws.Range("E98").Value2 = 0
$ws.Range("F98").Value2 = 0
$ws.Range("G98").Value2 = 0
$ws.Range("H98").Value2 = 0
$ws.Range("I98").Value2 = 0
$ws.Range("J98").Value2 = 0
$ws.Range("K98").Value2 = 0
$ws.Range("L98").Value2 = 0
$ws.Range("M98").Value2 = 0
$ws.Range("N98").Value2 = 0
# Row 99
$ws.Range("D99").Value2 = 0
$ws.Range("E99").Value2 = 0
$ws.Range("F99").Value2 = 0
$ws.Range("G99").Value2 = 0
$ws.Range("H99").Value2 = 0
$ws.Range("I99").Value2 = 0
$ws.Range("J99").Value2 = 0
$ws.Range("K99").Value2 = 0
$ws.Range("L99").Value2 = 0
$ws.Range("M99").Value2 = 0
$ws.Range("N99").Value2 = 0
# Row 100
$ws.Range("D100").Value2 = 41400
$ws.Range("E100").Value2 = 140300
$ws.Range("F100").Value2 = -18100
$ws.Range("G100").Value2 = 323000
$ws.Range("H100").Value2 = -21100
$ws.Range("I100").Value2 = -24400
$ws.Range("J100").Value2 = -271300
$ws.Range("K100").Value2 = -13900
$ws.Range("L100").Value2 = -15200
$ws.Range("M100").Value2 = -24000
$ws.Range("N100").Value2 = -20200
# Row 101
$ws.Range("D101").Value2 = 100
$ws.Range("E101").Value2 = -400
$ws.Range("F101").Value2 = -200
$ws.Range("G101").Value2 = -700
$ws.Range("H101").Value2 = 500
$ws.Range("I101").Value2 = -200
$ws.Range("J101").Value2 = -900
$ws.Range("K101").Value2 = 400
$ws.Range("L101").Value2 = 1000
$ws.Range("M101").Value2 = -1700
$ws.Range("N101").Value2 = -500
# Row 102
$ws.Range("D102").Value2 = 10600
$ws.Range("E102").Value2 = -575100
$ws.Range("F102").Value2 = -1100
$ws.Range("G102").Value2 = 405300
$ws.Range("H102").Value2 = 65000
$ws.Range("I102").Value2 = -122400
$ws.Range("J102").Value2 = -12100
$ws.Range("K102").Value2 = -119800
$ws.Range("L102").Value2 = -69600
$ws.Range("M102").Value2 = -52400
$ws.Range("N102").Value2 = 33700
Write-Host "done"
